$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 14 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2()
    $cell.Value = $old + 14
}

# Update the three changed production values (rows 41-43, column B)
$ws.Cells.Item(41, 2).Value = 759
$ws.Cells.Item(42, 2).Value = 789
$ws.Cells.Item(43, 2).Value = 846
